$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from an existing date cell (A2) onto the new
# row's date cell BEFORE assigning its value, so the new cell reuses the
# workbook's existing "s=1" (numFmtId 14) style instead of minting a new one.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$newDate = Get-Date -Year 2023 -Month 10 -Day 20
$ws.Range("A7").Value = $newDate.Date

$ws.Range("B7").Value = "10:00am"
$ws.Range("C7").Value = "10:45pm"
$ws.Range("D7").Value = "45 minutes"
$ws.Range("E7").Value = "Update and design concept review discussion"

# Widen column A (dates) and column E (main topic) to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 9.67
$ws.Columns.Item(5).ColumnWidth = 116.3

# Move the active selection the way the author's workbook ended up.
$ws.Range("E11").Select() | Out-Null
